$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.890575333333333
$ws.Range("H2").Value = 5.671726
$ws.Range("I2").Value = 0.006346320422088561
$ws.Range("J2").Value = 0.00634632042208856
$ws.Range("M2").Value = 20.56839166666667
$ws.Range("N2").Value = 61.705175
$ws.Range("O2").Value = 0.1304525281245593
$ws.Range("P2").Value = 0.1304525281245593
$ws.Range("Q2").Value = 38.88609393133889
$ws.Range("R2").Value = 349.97484538205
$ws.Range("S2").Value = 0.000827893543349973
$ws.Range("T2").Value = 0.0008278935433499731
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.890575333333333
$ws.Range("H3").Value = 5.671726
$ws.Range("I3").Value = 0.006346320422088561
$ws.Range("J3").Value = 0.00634632042208856
$ws.Range("O3").Value = 0.6526310778549473
$ws.Range("P3").Value = 0.6526310778549473
$ws.Range("Q3").Value = 194.5402956985749
$ws.Range("R3").Value = 1750.862661287174
$ws.Range("S3").Value = 0.004141805937480522
$ws.Range("T3").Value = 0.004141805937480521
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.890575333333333
$ws.Range("H4").Value = 5.671726
$ws.Range("I4").Value = 0.006346320422088561
$ws.Range("J4").Value = 0.00634632042208856
$ws.Range("O4").Value = 0.2169163940204933
$ws.Range("P4").Value = 0.2169163940204934
$ws.Range("Q4").Value = 64.65977619900355
$ws.Range("R4").Value = 581.9379857910319
$ws.Range("S4").Value = 0.001376620941258066
$ws.Range("T4").Value = 0.001376620941258066
# Row 5
$ws.Range("I5").Value = 0.8887896079640043
$ws.Range("J5").Value = 0.8887896079640044
$ws.Range("M5").Value = 20.56839166666667
$ws.Range("N5").Value = 61.705175
$ws.Range("O5").Value = 0.1304525281245593
$ws.Range("P5").Value = 0.1304525281245593
$ws.Range("Q5").Value = 5445.920451824903
$ws.Range("R5").Value = 49013.28406642412
$ws.Range("S5").Value = 0.1159448513297403
$ws.Range("T5").Value = 0.1159448513297403
# Row 6
$ws.Range("I6").Value = 0.8887896079640043
$ws.Range("J6").Value = 0.8887896079640044
$ws.Range("O6").Value = 0.6526310778549473
$ws.Range("P6").Value = 0.6526310778549473
$ws.Range("S6").Value = 0.5800517198318242
$ws.Range("T6").Value = 0.5800517198318242
# Row 7
$ws.Range("I7").Value = 0.8887896079640043
$ws.Range("J7").Value = 0.8887896079640044
$ws.Range("O7").Value = 0.2169163940204933
$ws.Range("P7").Value = 0.2169163940204934
$ws.Range("S7").Value = 0.1927930368024398
$ws.Range("T7").Value = 0.1927930368024398
# Row 8
$ws.Range("I8").Value = 0.104864071613907
$ws.Range("J8").Value = 0.104864071613907
$ws.Range("M8").Value = 20.56839166666667
$ws.Range("N8").Value = 61.705175
$ws.Range("O8").Value = 0.1304525281245593
$ws.Range("P8").Value = 0.1304525281245593
$ws.Range("Q8").Value = 642.5383320716501
$ws.Range("R8").Value = 5782.84498864485
$ws.Range("S8").Value = 0.013679783251469
$ws.Range("T8").Value = 0.013679783251469
# Row 9
$ws.Range("I9").Value = 0.104864071613907
$ws.Range("J9").Value = 0.104864071613907
$ws.Range("O9").Value = 0.6526310778549473
$ws.Range("P9").Value = 0.6526310778549473
$ws.Range("S9").Value = 0.06843755208564248
$ws.Range("T9").Value = 0.06843755208564248
# Row 10
$ws.Range("I10").Value = 0.104864071613907
$ws.Range("J10").Value = 0.104864071613907
$ws.Range("O10").Value = 0.2169163940204933
$ws.Range("P10").Value = 0.2169163940204934
$ws.Range("S10").Value = 0.02274673627679547
$ws.Range("T10").Value = 0.02274673627679548
